$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")
$ws.Activate()

# New card entry ("Rick Wise four cards") belongs right above the existing
# "Steve Carlton Traded" row, so insert a fresh row at 602 and push the rest
# of the list down by one.
$ws.Rows.Item(602).Insert()

$ws.Range("A602").Value = "Rick Wise four cards"
$ws.Range("B602").Value = "https://3.bp.blogspot.com/-gHwFpu8LpZk/WIkJAvGozGI/AAAAAAABr3U/kCL6Z-bpVY0uIWixiUwr9BohTm8jSbi7ACLcB/s1600/wise-2.jpg"

$ws.Range("B602").Select()
